# Update the "ecological" column p-values in the first table
# ("Cod: p-values for trend lines in Fig 8") to use Trange instead of T50.
#
# Table 1 layout (header row 1: Var | ecological | institutional | socioeconomic):
#   Row 2 GDP 2016      ecological 0.10 -> 0.07
#   Row 3 OHI economic  ecological 0.75 -> 0.80
#   Row 4 OHI fisheries ecological 0.51 -> 0.64
#   Row 5 Readiness     ecological 0.77 -> 0.85
#   Row 6 Vulnerability ecological 0.69 -> 0.82

$d = $word.ActiveDocument
$t = $d.Tables(1)

$changes = @(
    @{ Row = 2; Old = "0.10"; New = "0.07" },
    @{ Row = 3; Old = "0.75"; New = "0.80" },
    @{ Row = 4; Old = "0.51"; New = "0.64" },
    @{ Row = 5; Old = "0.77"; New = "0.85" },
    @{ Row = 6; Old = "0.69"; New = "0.82" }
)

# NB: Find.Execute on a Cell.Range object is not reliably scoped in this
# runtime -- it can match/replace occurrences elsewhere in the document
# (especially with ReplaceAll). Re-anchor a fresh Document.Range using the
# cell's Start/End offsets and use ReplaceOne so the edit stays confined to
# the single targeted cell (there are duplicate p-values, e.g. "0.69" and
# "0.75", elsewhere in these tables that must NOT change).
foreach ($chg in $changes) {
    $cellRange = $t.Cell($chg.Row, 2).Range
    $rng = $d.Range($cellRange.Start, $cellRange.End)
    $rng.Find.Execute($chg.Old, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $chg.New, 1) | Out-Null
}
